# Prepare presentation for lab meeting:
#  - Send the 6 colored "Rectangle" callout boxes to the back of the slide
#    (one at a time, in their original front-to-back order, so their
#    relative order reverses), and
#  - Swap their look from "no fill / colored outline" to "colored fill /
#    no outline".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName($slide, $name) {
    $n = $slide.Shapes.Count
    for ($i = 1; $i -le $n; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) { return $sh }
    }
    return $null
}

# VBA/COM ".RGB" values are stored little-endian (0xBBGGRR), so convert
# from a normal "RRGGBB" hex string.
function ConvertTo-RGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target fill color for each rectangle, keyed by its current shape name.
$fillColors = @{
    "Rectangle 4"   = "CAE3E4"
    "Rectangle 206" = "C4C5DE"
    "Rectangle 207" = "EDC4F8"
    "Rectangle 208" = "F8C4C4"
    "Rectangle 209" = "F6E0B8"
    "Rectangle 210" = "BCF5A5"
}

# Send to back in this order so the final stacking order (front-to-back in
# the Shapes collection / back-to-front in the XML) ends up reversed:
# Rectangle 210, 209, 208, 207, 206, 4.
$sendOrder = @("Rectangle 4", "Rectangle 206", "Rectangle 207", "Rectangle 208", "Rectangle 209", "Rectangle 210")

foreach ($nm in $sendOrder) {
    $sh = Get-ShapeByName $s $nm
    $sh.ZOrder(1)
}

foreach ($nm in $fillColors.Keys) {
    $sh = Get-ShapeByName $s $nm
    $sh.Fill.Visible = -1
    $sh.Fill.ForeColor.RGB = ConvertTo-RGB $fillColors[$nm]
    $sh.Line.Visible = 0
}
